# Update cell values in Sheet1 to reflect the updated RandomForest imputation
# results (commit: "Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B4" = 8.6457
    "C4" = -12.10579999999999
    "B7" = 5.381200000000002
    "D10" = -8.000099999999996
    "C12" = -10.4707
    "D13" = -8.355899999999993
    "B16" = 4.8792
    "C18" = -11.9074
    "C19" = -11.2969
    "C20" = -12.9354
    "B28" = 6.357699999999999
    "B29" = 5.150900000000003
    "D30" = -7.022
    "C31" = -13.5991
    "B32" = 6.695799999999999
    "B40" = 9.167499999999995
    "C40" = -13.14990000000001
    "D40" = -9.165199999999995
    "C42" = -11.46880000000001
    "D44" = -6.709100000000004
    "C47" = -12.06189999999999
    "C48" = -11.8378
    "B52" = 5.233799999999998
    "B57" = 4.755999999999995
    "C63" = -10.3054
    "C64" = -10.85989999999999
    "B66" = 5.429399999999998
    "C76" = -12.06350000000001
    "C81" = -14.3334
    "C89" = -13.19609999999999
    "D89" = -8.275999999999994
    "D91" = -8.16
    "C94" = -10.77589999999999
    "B100" = 4.859000000000002
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
